$d = $word.ActiveDocument

# 1. Merge the split "крита" run (with spell-check proofErr wrapper) back into
#    a single contiguous phrase.
$d.Content.Find.Execute("шанс уклона, шанс крита, шанс блока", $true, $false, $false, $false, $false,
                         $true, 1, $false, "шанс уклона, шанс крита, шанс блока", 2) | Out-Null

# 2. Update the "Крит" description text.
$d.Content.Find.Execute("процентное увеличение наносимого урона, усиливает накладываемые эффекты", $true, $false, $false, $false, $false,
                         $true, 1, $false, "процентный шанс крита, усиливает накладываемые эффекты", 2) | Out-Null

# 3. Update the "Сила" description text.
$d.Content.Find.Execute("урон, процент увеличения урона при критическом ударе", $true, $false, $false, $false, $false,
                         $true, 1, $false, "бонусный урон, процент увеличения урона при критическом ударе", 2) | Out-Null
